$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
Write-Host $ws.Name
$ws.Range("H2").Value2 = 767.8
